$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing last row (215): D and F changed ---
$ws.Range("D215").Value = 1309
$ws.Range("F215").Value = 1308.62

# --- Append new rows 216-218, cloning row 215's formatting first so the
#     date column keeps its existing date-number-format style (s="2") ---
$ws.Range("A215:G215").Copy($ws.Range("A216:G216"))
$ws.Range("A215:G215").Copy($ws.Range("A217:G217"))
$ws.Range("A215:G215").Copy($ws.Range("A218:G218"))

# Row 216
$ws.Range("A216").Value = 45170.33333333334
$ws.Range("B216").Value = "FX_IDC:USDIQD"
$ws.Range("C216").Value = 1308.62
$ws.Range("D216").Value = 1308.65
$ws.Range("E216").Value = 1307
$ws.Range("F216").Value = 1308.62
$ws.Range("G216").Value = 0

# Row 217
$ws.Range("A217").Value = 45201.375
$ws.Range("B217").Value = "FX_IDC:USDIQD"
$ws.Range("C217").Value = 1308.62
$ws.Range("D217").Value = 1310
$ws.Range("E217").Value = 1306
$ws.Range("F217").Value = 1309
$ws.Range("G217").Value = 0

# Row 218
$ws.Range("A218").Value = 45231.375
$ws.Range("B218").Value = "FX_IDC:USDIQD"
$ws.Range("C218").Value = 1309
$ws.Range("D218").Value = 1309
$ws.Range("E218").Value = 1308
$ws.Range("F218").Value = 1309
$ws.Range("G218").Value = 0
